$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) sometimes holds plain decimal-looking text
# (e.g. "351.98"). Excel would normally auto-convert such text into a
# real number, losing the original text formatting/trailing zeros, so
# we explicitly force the Text number format for those specific cells
# before writing the value, matching the source data which stores all
# of column D and E as text.

$ws.Range("D2").Value = "51.706.33"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "2.812.24"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "351.98"
$ws.Range("E5").Value = "  +5.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "112.74"
$ws.Range("E6").Value = "  -3.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.565"
$ws.Range("E7").Value = "  +4.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.598"
$ws.Range("E9").Value = "  +3.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.37"
$ws.Range("E10").Value = "  -1.73%  "
$ws.Range("E11").Value = "  -1.46%  "
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("E13").Value = "  -2.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.72"
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("D15").Value = "3.254.67"
$ws.Range("E15").Value = "  +1.55%  "
$ws.Range("D16").Value = "2.812.11"
$ws.Range("E16").Value = "  +1.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.883"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").Value = "51.472.61"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.47"
$ws.Range("E19").Value = "  +8.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.19"
$ws.Range("E21").Value = "  -1.63%  "
$ws.Range("D22").Value = "0.0₃0991"
$ws.Range("E22").Value = "  +1.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "270.41"
$ws.Range("E23").Value = "  -2.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.58"
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.73"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.64"
$ws.Range("E26").Value = "  -0.74%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.30"
$ws.Range("E28").Value = "  +1.10%  "
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.139"
$ws.Range("E30").Value = "  -2.20%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.94"
$ws.Range("E31").Value = "  -3.45%  "
$ws.Range("B32").Value = "OKB"
$ws.Range("C32").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.53"
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("B33").Value = "VeChain"
$ws.Range("C33").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0447"
$ws.Range("E33").Value = "  +25.53%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.80"
$ws.Range("E34").Value = "  +3.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0822"
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.16"
$ws.Range("E36").Value = "  +2.49%  "
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.06"
$ws.Range("E38").Value = "  -1.84%  "
$ws.Range("E39").Value = "  -1.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.12"
$ws.Range("E40").Value = "  -6.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.64"
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.117"
$ws.Range("E42").Value = "  +1.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "126.10"
$ws.Range("E43").Value = "  -1.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.52"
$ws.Range("E44").Value = "  +2.41%  "
$ws.Range("E45").Value = "  -1.01%  "
$ws.Range("D46").Value = "2.075.63"
$ws.Range("E46").Value = "  -0.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.32"
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.28"
$ws.Range("E48").Value = "  +2.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.64"
$ws.Range("E49").Value = "  +1.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.929"
$ws.Range("E50").Value = "  +5.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "60.66"
$ws.Range("E51").Value = "  +0.33%  "
